$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45203 -> 45204) for every data row (rows 2 through 163).
for ($r = 2; $r -le 163; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value = 45204
    }
}
